# "System of life has been implemented. Game for restoring of life is in progress."
#
# Adds the next four story beats (rows 16-19) to the Events sheet, each a
# (id, story) pair appended under the existing rows, and moves the active
# selection to B8 (where the author was working next).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = 30
$ws.Range("B16").Value = @'
Мы шли по шоссе целый день, попутно говоря о данной ситуации, зомби встречались, но встречались очень редко и обойти их или убить не составляло труда. Уже начало смеркаться. Мы думали о ночлеге, так как за день дороги очень устали. Но вдруг вдалеке зажглись огни. Несомненно, это был военный блокпост. Это было понятно, потому что были видны четыре вышки, которые светили прожекторами на прилежавшую территорию.
'@

$ws.Range("A17").Value = 41
$ws.Range("B17").Value = @'
Мы расставили платки недалеко от шоссе, так как боялись потеряться. Мы решили, что первый на стрёме стою я. Прошло уже 3 часа. Уже была очередь сменяться. Но из-за своей сонливости я не заметил, как зомби подошёл ко мне и укусил. Чувствовалась огромная боль в моей руке. Мне удалось убить его, но на мои крики стекались зомби из леса. Пока Виктор вставал на меня напали уже двое ходячих. Мы отбились и от них. Но нас уже окружили около дюжины зомби. И сначала меня, а потом и Виктора повалили и загрызли. Смерть 
'@

$ws.Range("A18").Value = 42
$ws.Range("B18").Value = @'
Мы сразу забыли о ночлеге и решили марш броском пройти это расстояние. Но приближаясь к блокпосту количество, зомби начало увеличиваться. Поначалу мы убивали их с помощью охотничьих ножей, но скоро стало совсем темно. Мы не заметили, как нас окружила толпа из примерно 30 зомби, благо они поначалу не замечали нас и время на то, чтобы придумать план действий был.
'@

$ws.Range("A19").Value = 52
$ws.Range("B19").Value = @'
Мы решили просто пробежать их. Поначалу у нас всё складывалось отлично, но из-за кромешной тьмы Виктор не увидел зомби и запнулся об него, у которого не было ног. Помочь ему мне сразу не удалось, так как мне загородили проход к нему три зомби.
'@

# Author's cursor ends up on B8 (next story cell to revisit) and the workbook
# window is scrolled down slightly - mirrored here via the selection + window
# position (the latter is cosmetic/session state).
$ws.Range("B8").Select() | Out-Null
$excel.ActiveWindow.Top = 1905

Write-Output "done"
